# Slide 1, shape "TextBox 5" ("Presented By: ... 1. Rajamanivarma Information Technology")
#
# The 2nd paragraph's trailing run " Information Technology" (which sits right
# after the "Rajamanivarma" run) is split into three runs:
#   "- "
#   "st.joseph"
#   " college of engineering - Information Technology"
#
# Net effect: "1. Rajamanivarma Information Technology"
#          -> "1. Rajamanivarma- st.joseph college of engineering - Information Technology"
#
# The shape uses spAutoFit, so PowerPoint grows the textbox height to fit the
# extra wrapped line automatically once the text changes.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)            # "TextBox 5"
$tr = $sh.TextFrame.TextRange

# 2nd paragraph: "1. Rajamanivarma Information Technology"
$para2 = $tr.Paragraphs(2, 1)

# 3rd run of that paragraph is " Information Technology" (1-based: "1. " / "Rajamanivarma" / " Information Technology")
$tailRun = $para2.Runs(3, 1)

# Re-purpose that run for the new "- " text, keeping its original formatting,
# then append the two remaining chunks as new runs right after it (inheriting
# formatting from their predecessor, matching the target rPr).
$tailRun.Text = "- "
$run2 = $tailRun.InsertAfter("st.joseph")
$run3 = $run2.InsertAfter(" college of engineering - Information Technology")
